# Scheduled runner update: refresh currentAveragePrice/Leve profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with latest market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 501.4
$ws.Range("I53").Value = 805.8333
$ws.Range("K53").Value = 805.8333
$ws.Range("M53").Value = -168.8333

$ws.Range("H129").Value = 2791.5881
$ws.Range("I129").Value = 1063
$ws.Range("K129").Value = 3189
$ws.Range("M129").Value = 1811

$ws.Range("H137").Value = 1520
$ws.Range("I137").Value = 1787.4286
$ws.Range("K137").Value = 5362.2858
$ws.Range("M137").Value = -2812.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1828.7142
$ws.Range("I45").Value = 1828.7142
$ws.Range("K45").Value = 1828.7142
$ws.Range("M45").Value = -1451.7142

$ws.Range("H61").Value = 1898.75
$ws.Range("I61").Value = 1898.75
$ws.Range("K61").Value = 1898.75
$ws.Range("M61").Value = -1686.75

$ws.Range("H63").Value = 2553.7778
$ws.Range("J63").Value = 3221.2
$ws.Range("L63").Value = 3221.2
$ws.Range("N63").Value = -4593.2

$ws.Range("H66").Value = 2553.7778
$ws.Range("J66").Value = 3221.2
$ws.Range("L66").Value = 16106
$ws.Range("N66").Value = -22970

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H136").Value = 1898.75
$ws.Range("I136").Value = 1898.75
$ws.Range("K136").Value = 5696.25
$ws.Range("M136").Value = -3146.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25216.889
$ws.Range("I82").Value = 6738
$ws.Range("K82").Value = 6738
$ws.Range("M82").Value = -6355

$ws.Range("H85").Value = 25216.889
$ws.Range("I85").Value = 6738
$ws.Range("K85").Value = 6738
$ws.Range("M85").Value = -5412

$ws.Range("H94").Value = 1568.1428
$ws.Range("J94").Value = 983.7143
$ws.Range("L94").Value = 983.7143
$ws.Range("N94").Value = -1885.7143

$ws.Range("H134").Value = 10017.818
$ws.Range("I134").Value = 10899.5
$ws.Range("K134").Value = 32698.5
$ws.Range("M134").Value = -30163.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 78.29412000000001
$ws.Range("I7").Value = 96.09999999999999
$ws.Range("J7").Value = 52.857143
$ws.Range("K7").Value = 96.09999999999999
$ws.Range("L7").Value = 52.857143
$ws.Range("M7").Value = 16.90000000000001
$ws.Range("N7").Value = -278.857143

$ws.Range("H22").Value = 244.88889
$ws.Range("I22").Value = 217
$ws.Range("J22").Value = 300.66666
$ws.Range("K22").Value = 217
$ws.Range("L22").Value = 300.66666
$ws.Range("M22").Value = 133
$ws.Range("N22").Value = -1000.66666

$ws.Range("H58").Value = 2265.1667
$ws.Range("I58").Value = 2219.4
$ws.Range("J58").Value = 2494
$ws.Range("K58").Value = 2219.4
$ws.Range("L58").Value = 2494
$ws.Range("M58").Value = -2016.4
$ws.Range("N58").Value = -2900

$ws.Range("H60").Value = 25000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 25000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 25000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -26022

$ws.Range("H94").Value = 1718.5555
$ws.Range("I94").Value = 1861.5714
$ws.Range("J94").Value = 1218
$ws.Range("K94").Value = 1861.5714
$ws.Range("L94").Value = 1218
$ws.Range("M94").Value = -1410.5714
$ws.Range("N94").Value = -2120

$ws.Range("H136").Value = 2265.1667
$ws.Range("I136").Value = 2219.4
$ws.Range("J136").Value = 2494
$ws.Range("K136").Value = 6658.200000000001
$ws.Range("L136").Value = 7482
$ws.Range("M136").Value = -4108.200000000001
$ws.Range("N136").Value = -12582

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3900
$ws.Range("J131").Value = 3900
$ws.Range("L131").Value = 11700
$ws.Range("N131").Value = -21780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 27960
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 27960
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H122").Value = 8933687
$ws.Range("I122").Value = 8933687
$ws.Range("K122").Value = 26801061
$ws.Range("M122").Value = -26798611

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1668.625
$ws.Range("J22").Value = 1749.8334
$ws.Range("L22").Value = 1749.8334
$ws.Range("N22").Value = -2339.8334

$ws.Range("H27").Value = 1668.625
$ws.Range("J27").Value = 1749.8334
$ws.Range("L27").Value = 1749.8334
$ws.Range("N27").Value = -1963.8334

$ws.Range("H57").Value = 6500
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 8000
$ws.Range("K57").Value = 5000
$ws.Range("L57").Value = 8000
$ws.Range("M57").Value = -4434
$ws.Range("N57").Value = -9132

$ws.Range("H93").Value = 386.75
$ws.Range("I93").Value = 399
$ws.Range("J93").Value = 350
$ws.Range("K93").Value = 399
$ws.Range("L93").Value = 350
$ws.Range("M93").Value = 849
$ws.Range("N93").Value = -2846

$ws.Range("H122").Value = 3290.5
$ws.Range("I122").Value = 3290.5
$ws.Range("K122").Value = 9871.5
$ws.Range("M122").Value = -7421.5

$ws.Range("H132").Value = 3091
$ws.Range("I132").Value = 3091
$ws.Range("K132").Value = 9273
$ws.Range("M132").Value = -6743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 237.5
$ws.Range("I81").Value = 300
$ws.Range("K81").Value = 600
$ws.Range("M81").Value = 461

$ws.Range("H84").Value = 237.5
$ws.Range("I84").Value = 300
$ws.Range("K84").Value = 3000
$ws.Range("M84").Value = 2304

$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 2076.75
$ws.Range("I136").Value = 2073.4285
$ws.Range("J136").Value = 2100
$ws.Range("K136").Value = 6220.2855
$ws.Range("L136").Value = 6300
$ws.Range("M136").Value = -3670.2855
$ws.Range("N136").Value = -11400
